## Fixed variables and query errors in Bread from TC01 to TC30
## This targets the "startup" sheet's CasesTab query cell (B2), which
## incorrectly returned an extra `Cohort` column that isn't part of the
## Cases result set. Remove that trailing RETURN line from the query text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$casesQuery = @'
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)
WHERE demo.breed IN ['Basset Hound']
MATCH (c)<--(diag:diagnosis)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (co:cohort)<-[*]-(c)
WITH DISTINCT c, s, demo, diag, co
RETURN  coalesce(c.case_id, '') AS `Case ID` ,
        coalesce(s.clinical_study_designation, '') AS `Study Code` ,
        coalesce(s.clinical_study_type, '') AS  `Study Type`,
        coalesce(demo.breed, '') AS Breed ,
        coalesce(diag.disease_term, '') AS Diagnosis ,
        coalesce(diag.stage_of_disease, '') AS `Stage of Disease` ,
        coalesce(demo.patient_age_at_enrollment, '') AS Age ,
        coalesce(demo.sex, '') AS Sex ,
        coalesce(demo.neutered_indicator, '') AS `Neutered Status`,
        coalesce(demo.weight, '') AS `Weight (kg)`,
        coalesce(diag.best_response, '') AS `Response to Treatment`

'@

# Row 2 (CasesTab) query text: drop the trailing `Cohort` column line - it
# does not exist in the CasesTab result set and was causing a query error.
$ws.Range("B2").Value = $casesQuery

# The view had scrolled/zoomed out while debugging; restore a normal
# working zoom level, with the cursor back on the fixed cell.
$win = $excel.ActiveWindow
$win.Zoom = 130
$ws.Range("B2").Select() | Out-Null

# Row heights shrink slightly to match the content/font re-measurement that
# happened when the workbook was re-saved (one fewer wrapped line in B2,
# plus a small global re-layout of the other wrapped rows).
$ws.Rows.Item(2).RowHeight = 259.2
$ws.Rows.Item(3).RowHeight = 230.4
$ws.Rows.Item(4).RowHeight = 244.8
